$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values for rows 3-7
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 100
$ws1.Range("F4").Value = 89
$ws1.Range("F5").Value = 2607
$ws1.Range("F6").Value = 243
$ws1.Range("F7").Value = 382

# Sheet "全部类型" (sheet4): update corresponding values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 100
$ws4.Range("F4").Value = 89
$ws4.Range("F5").Value = 2607
$ws4.Range("F6").Value = 243
$ws4.Range("F9").Value = 382
